$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sciences")

# --- 1) Push the bottom SUBTOTAL row (68) down to row 69, preserving its styles ---
$ws.Range("A68:B68").Copy($ws.Range("A69:B69"))
$ws.Range("A69").Formula = "=SUBTOTAL(2,A2:A68)"
$ws.Range("B69").Formula = "=SUBTOTAL(9,B2:B68)"

# --- 2) Duplicate the (still original) row 28 data into the freed-up row 68, preserving styles ---
$ws.Range("A28:M28").Copy($ws.Range("A68:M68"))
$ws.Range("A68").Value2 = 67
$ws.Range("B68").Formula = "=IF(OR(NOT(ISERROR(SEARCH(""archive.org"",M68))),NOT(ISERROR(SEARCH(""app.box.com"",M68))),NOT(ISERROR(SEARCH(""islamway.net"",M68))),NOT(ISERROR(SEARCH(""qurancomplex.gov.sa"",M68))),NOT(ISERROR(SEARCH(""tanzil.net"",M68))),NOT(ISERROR(SEARCH(""alsirah.com"",M68))),NOT(ISERROR(SEARCH(""i36"",M68))),(RIGHT(M68,4)="".pdf""),C68=6,C68=8,C68=9),0,1)"

Write-Output "step1-2 done"
